$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.928.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.975.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.10%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.24%  "

$ws.Range("E6").Value = "  -3.78%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.369"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.77%  "

$ws.Range("E11").Value = "  -6.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0979"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.261.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.747"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.76%  "

$ws.Range("E17").Value = "  -6.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.984.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.791.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0803"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.37%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.67%  "

$ws.Range("E32").Value = "  -3.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0609"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.39%  "

$ws.Range("E35").Value = "  -6.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.67%  "

$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.416.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("E43").Value = "  -5.88%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0202"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.73%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0885"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.74%  "

$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.46%  "

Write-Output "applied"